# Append 20 new daily-position rows (284-303) to the bottom of the sheet,
# covering 2025-02-26 .. 2025-03-17, matching the columns that already
# carry data (B, C, I, K, N, O, Q, U, Z); the remaining columns for those
# dates have no quotes and stay blank, same as the existing rows above.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A284").Value = 45714
$ws.Range("B284").Value = 745.6790115702
$ws.Range("C284").Value = 161.9953680495
$ws.Range("I284").Value = 242.282760292
$ws.Range("K284").Value = 110.280232945056
$ws.Range("N284").Value = 19.208012208
$ws.Range("O284").Value = 1.0411795779
$ws.Range("Q284").Value = 0.0000019992
$ws.Range("U284").Value = 238.2038317619702
$ws.Range("Z284").Value = 3313.552671284544

$ws.Range("A285").Value = 45715
$ws.Range("B285").Value = 749.7370056924001
$ws.Range("C285").Value = 160.008881622
$ws.Range("I285").Value = 246.508206008
$ws.Range("K285").Value = 110.959776633978
$ws.Range("N285").Value = 19.4408365984
$ws.Range("O285").Value = 1.0306923627
$ws.Range("Q285").Value = 0.0000019704
$ws.Range("U285").Value = 235.7731804174603
$ws.Range("Z285").Value = 3458.744187203262

$ws.Range("A286").Value = 45716
$ws.Range("B286").Value = 746.5627619532
$ws.Range("C286").Value = 155.1463233965
$ws.Range("I286").Value = 265.307858558
$ws.Range("K286").Value = 113.241101875359
$ws.Range("N286").Value = 19.39427172032
$ws.Range("O286").Value = 1.0011715329
$ws.Range("Q286").Value = 0.0000019008
$ws.Range("U286").Value = 233.7263161273467
$ws.Range("Z286").Value = 3511.799412301568

$ws.Range("A287").Value = 45717
$ws.Range("B287").Value = 761.7382208334001
$ws.Range("C287").Value = 153.7457291265
$ws.Range("I287").Value = 257.215055746
$ws.Range("K287").Value = 113.629412554743
$ws.Range("N287").Value = 18.26507342688
$ws.Range("O287").Value = 1.0349825871
$ws.Range("Q287").Value = 0.0000018216
$ws.Range("U287").Value = 228.3532973657985
$ws.Range("Z287").Value = 3536.906318647474

$ws.Range("A288").Value = 45718
$ws.Range("B288").Value = 834.3630306000001
$ws.Range("C288").Value = 174.5965562985
$ws.Range("I288").Value = 319.970086401
$ws.Range("K288").Value = 118.483296047043
$ws.Range("N288").Value = 20.39541659904
$ws.Range("O288").Value = 1.0622391318
$ws.Range("Q288").Value = 0.0000021192
$ws.Range("U288").Value = 259.5679777900309
$ws.Range("Z288").Value = 3560.243164805728

$ws.Range("A289").Value = 45719
$ws.Range("B289").Value = 763.1196505758001
$ws.Range("C289").Value = 149.0045095135
$ws.Range("I289").Value = 254.529391096
$ws.Range("K289").Value = 112.221786341976
$ws.Range("N289").Value = 16.92633318208
$ws.Range("O289").Value = 0.9807929670000001
$ws.Range("Q289").Value = 0.0000017352
$ws.Range("U289").Value = 218.8865500240231
$ws.Range("Z289").Value = 3336.237390005242

$ws.Range("A290").Value = 45720
$ws.Range("B290").Value = 772.5136029444001
$ws.Range("C290").Value = 150.5645773885
$ws.Range("I290").Value = 259.238256449
$ws.Range("K290").Value = 117.755213523198
$ws.Range("N290").Value = 15.5410280592
$ws.Range("O290").Value = 0.9942595047
$ws.Range("Q290").Value = 0.0000016776
$ws.Range("U290").Value = 215.4324665344564
$ws.Range("Z290").Value = 3341.500990036944

$ws.Range("A291").Value = 45721
$ws.Range("B291").Value = 801.9338611878001
$ws.Range("C291").Value = 155.4236687965
$ws.Range("I291").Value = 261.906016668
$ws.Range("K291").Value = 118.289140707351
$ws.Range("N291").Value = 16.8215622064
$ws.Range("O291").Value = 1.0184856528
$ws.Range("Q291").Value = 0.0000016992
$ws.Range("U291").Value = 226.9460781663454
$ws.Range("Z291").Value = 3307.264299565254

$ws.Range("A292").Value = 45722
$ws.Range("B292").Value = 795.9673733742001
$ws.Range("C292").Value = 152.69250997
$ws.Range("I292").Value = 256.588400661
$ws.Range("K292").Value = 116.541742650123
$ws.Range("N292").Value = 16.32098976704
$ws.Range("O292").Value = 1.015523355
$ws.Range("Q292").Value = 0.0000016392
$ws.Range("U292").Value = 218.5027629696268
$ws.Range("Z292").Value = 3179.866546585564

$ws.Range("A293").Value = 45723
$ws.Range("B293").Value = 768.2631928650001
$ws.Range("C293").Value = 148.49072716
$ws.Range("I293").Value = 249.569863709
$ws.Range("K293").Value = 118.240601872428
$ws.Range("N293").Value = 15.66908147392
$ws.Range("O293").Value = 1.0125780819
$ws.Range("Q293").Value = 0.0000016752
$ws.Range("U293").Value = 210.3153058091724
$ws.Range("Z293").Value = 3062.157544106706

$ws.Range("A294").Value = 45724
$ws.Range("B294").Value = 763.135936011
$ws.Range("C294").Value = 152.788194133
$ws.Range("I294").Value = 245.308609131
$ws.Range("K294").Value = 118.337679542274
$ws.Range("N294").Value = 15.22671513216
$ws.Range("O294").Value = 1.009734957
$ws.Range("Q294").Value = 0.0000016008
$ws.Range("U294").Value = 206.8612223196057
$ws.Range("Z294").Value = 3045.248811261504

$ws.Range("A295").Value = 45725
$ws.Range("B295").Value = 714.5621473086001
$ws.Range("C295").Value = 140.0878549035
$ws.Range("I295").Value = 226.508956581
$ws.Range("K295").Value = 112.658635856283
$ws.Range("N295").Value = 13.10801317952
$ws.Range("O295").Value = 0.9442749855
$ws.Range("Q295").Value = 0.0000014088
$ws.Range("U295").Value = 174.2393226959202
$ws.Range("Z295").Value = 2876.999932372056

$ws.Range("A296").Value = 45726
$ws.Range("B296").Value = 695.6346657708001
$ws.Range("C296").Value = 129.319226385
$ws.Range("I296").Value = 211.845227592
$ws.Range("K296").Value = 111.251009643516
$ws.Range("N296").Value = 12.47938732544
$ws.Range("O296").Value = 0.9048457803
$ws.Range("Q296").Value = 0.000001368
$ws.Range("U296").Value = 161.0626338283139
$ws.Range("Z296").Value = 2929.542771626568

$ws.Range("A297").Value = 45727
$ws.Range("B297").Value = 734.0216492322002
$ws.Range("C297").Value = 133.3636156805
$ws.Range("I297").Value = 224.432042585
$ws.Range("K297").Value = 108.921145567212
$ws.Range("N297").Value = 12.96831854528
$ws.Range("O297").Value = 0.9390994767
$ws.Range("Q297").Value = 0.0000015192
$ws.Range("U297").Value = 167.9708008074473
$ws.Range("Z297").Value = 2940.302874346242

$ws.Range("A298").Value = 45728
$ws.Range("B298").Value = 740.6343324936
$ws.Range("C298").Value = 132.30762307
$ws.Range("I298").Value = 226.705905322
$ws.Range("K298").Value = 108.193063043367
$ws.Range("N298").Value = 14.01602830208
$ws.Range("O298").Value = 0.9701865789
$ws.Range("Q298").Value = 0.0000016704
$ws.Range("U298").Value = 174.1113936777881
$ws.Range("Z298").Value = 3039.286503260992

$ws.Range("A299").Value = 45729
$ws.Range("B299").Value = 717.9379233084001
$ws.Range("C299").Value = 129.2838648465
$ws.Range("I299").Value = 220.886965247
$ws.Range("K299").Value = 109.357995081519
$ws.Range("N299").Value = 13.63186805792
$ws.Range("O299").Value = 0.9870921059999999
$ws.Range("Q299").Value = 0.0000016056
$ws.Range("U299").Value = 165.028433390409
$ws.Range("Z299").Value = 2941.32764603383

$ws.Range("A300").Value = 45730
$ws.Range("B300").Value = 743.3168268960001
$ws.Range("C300").Value = 132.5468334775
$ws.Range("I300").Value = 239.095771574
$ws.Range("K300").Value = 107.950368868752
$ws.Range("N300").Value = 14.09751683872
$ws.Range("O300").Value = 1.0006437672
$ws.Range("Q300").Value = 0.0000016944
$ws.Range("U300").Value = 175.262754840977
$ws.Range("Z300").Value = 2914.124615781494

$ws.Range("A301").Value = 45731
$ws.Range("B301").Value = 746.4609779832001
$ws.Range("C301").Value = 134.3162971295
$ws.Range("I301").Value = 243.249599566
$ws.Range("K301").Value = 107.513519354445
$ws.Range("N301").Value = 14.3187000096
$ws.Range("O301").Value = 1.0552419801
$ws.Range("Q301").Value = 0.0000017136
$ws.Range("U301").Value = 187.9277276360549
$ws.Range("Z301").Value = 3084.749101764896

$ws.Range("A302").Value = 45732
$ws.Range("B302").Value = 730.8489986334
$ws.Range("C302").Value = 130.83769245
$ws.Range("I302").Value = 225.828588203
$ws.Range("K302").Value = 102.853791201837
$ws.Range("N302").Value = 13.71335659456
$ws.Range("O302").Value = 1.025227434
$ws.Range("Q302").Value = 0.0000015792
$ws.Range("U302").Value = 173.7276066233918
$ws.Range("Z302").Value = 3148.005463207828

$ws.Range("A303").Value = 45733
$ws.Range("B303").Value = 743.5542933234001
$ws.Range("C303").Value = 133.5633043685
$ws.Range("I303").Value = 229.194621231
$ws.Range("K303").Value = 107.076669840138
$ws.Range("N303").Value = 14.81927244896
$ws.Range("O303").Value = 1.0754332743
$ws.Range("Q303").Value = 0.0000017736
$ws.Range("U303").Value = 180.5078445843931
$ws.Range("Z303").Value = 3260.916670967524

# New dates in column A should keep the same date-stamp style (border,
# bold font, centered, yyyy-mm-dd hh:mm:ss format) as the rows above —
# copy formats only so we reuse the existing style record instead of
# minting a duplicate one.
$ws.Range("A283").Copy()
$ws.Range("A284:A303").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
